$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = "domain"
